$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).NumberFormat = "@"
$ws.Cells.Item(2,4).Value = '30.815.60'
$ws.Cells.Item(2,4).Style = "Normal"
$ws.Cells.Item(2,5).NumberFormat = "@"
$ws.Cells.Item(2,5).Value = '  +2.56%  '
$ws.Cells.Item(2,5).Style = "Normal"

$ws.Cells.Item(3,4).NumberFormat = "@"
$ws.Cells.Item(3,4).Value = '1.899.44'
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Cells.Item(3,5).NumberFormat = "@"
$ws.Cells.Item(3,5).Value = '  +0.74%  '
$ws.Cells.Item(3,5).Style = "Normal"

$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = '1.003'
$ws.Cells.Item(4,4).Style = "Normal"
$ws.Cells.Item(4,5).NumberFormat = "@"
$ws.Cells.Item(4,5).Value = '  +0.25%  '
$ws.Cells.Item(4,5).Style = "Normal"

$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = '245.94'
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).NumberFormat = "@"
$ws.Cells.Item(5,5).Value = '  +1.29%  '
$ws.Cells.Item(5,5).Style = "Normal"

$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = '1.002'
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).NumberFormat = "@"
$ws.Cells.Item(6,5).Value = '  +0.13%  '
$ws.Cells.Item(6,5).Style = "Normal"

$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = '0.4979'
$ws.Cells.Item(7,4).Style = "Normal"
$ws.Cells.Item(7,5).NumberFormat = "@"
$ws.Cells.Item(7,5).Value = '  +0.40%  '
$ws.Cells.Item(7,5).Style = "Normal"

$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = '0.2979'
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,5).NumberFormat = "@"
$ws.Cells.Item(8,5).Value = '  +1.73%  '
$ws.Cells.Item(8,5).Style = "Normal"

$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = '0.06832'
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).NumberFormat = "@"
$ws.Cells.Item(9,5).Value = '  +3.46%  '
$ws.Cells.Item(9,5).Style = "Normal"

$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = '1.908.53'
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).NumberFormat = "@"
$ws.Cells.Item(10,5).Value = '  +1.28%  '
$ws.Cells.Item(10,5).Style = "Normal"

$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = '17.27'
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).NumberFormat = "@"
$ws.Cells.Item(11,5).Value = '  +3.07%  '
$ws.Cells.Item(11,5).Style = "Normal"

$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = '0.07331'
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,5).NumberFormat = "@"
$ws.Cells.Item(12,5).Value = '  +2.26%  '
$ws.Cells.Item(12,5).Style = "Normal"

$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = '91.82'
$ws.Cells.Item(13,4).Style = "Normal"
$ws.Cells.Item(13,5).NumberFormat = "@"
$ws.Cells.Item(13,5).Value = '  +6.93%  '
$ws.Cells.Item(13,5).Style = "Normal"

$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = '5.097'
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).NumberFormat = "@"
$ws.Cells.Item(14,5).Value = '  +5.36%  '
$ws.Cells.Item(14,5).Style = "Normal"

$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = '0.6789'
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(15,5).NumberFormat = "@"
$ws.Cells.Item(15,5).Value = '  +2.45%  '
$ws.Cells.Item(15,5).Style = "Normal"

$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = '30.866.58'
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Cells.Item(16,5).NumberFormat = "@"
$ws.Cells.Item(16,5).Value = '  +2.86%  '
$ws.Cells.Item(16,5).Style = "Normal"

$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = '0.000008067'
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Cells.Item(17,5).NumberFormat = "@"
$ws.Cells.Item(17,5).Value = '  +1.45%  '
$ws.Cells.Item(17,5).Style = "Normal"

$ws.Cells.Item(18,5).NumberFormat = "@"
$ws.Cells.Item(18,5).Value = '  +5.58%  '
$ws.Cells.Item(18,5).Style = "Normal"

$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = '0.9999'
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(19,5).NumberFormat = "@"
$ws.Cells.Item(19,5).Value = '  -0.04%  '
$ws.Cells.Item(19,5).Style = "Normal"

$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = '2.153.91'
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).NumberFormat = "@"
$ws.Cells.Item(20,5).Value = '  +1.49%  '
$ws.Cells.Item(20,5).Style = "Normal"

$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = '1.000'
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).NumberFormat = "@"
$ws.Cells.Item(21,5).Value = '  +0.10%  '
$ws.Cells.Item(21,5).Style = "Normal"

$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = '4.878'
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).NumberFormat = "@"
$ws.Cells.Item(22,5).Value = '  +2.66%  '
$ws.Cells.Item(22,5).Style = "Normal"

$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = '181.99'
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(23,5).NumberFormat = "@"
$ws.Cells.Item(23,5).Value = '  +34.76%  '
$ws.Cells.Item(23,5).Style = "Normal"

$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = '6.096'
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(24,5).NumberFormat = "@"
$ws.Cells.Item(24,5).Value = '  +9.09%  '
$ws.Cells.Item(24,5).Style = "Normal"

$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = '9.362'
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,5).NumberFormat = "@"
$ws.Cells.Item(25,5).Value = '  +2.99%  '
$ws.Cells.Item(25,5).Style = "Normal"

$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = '154.79'
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Cells.Item(26,5).NumberFormat = "@"
$ws.Cells.Item(26,5).Value = '  +3.17%  '
$ws.Cells.Item(26,5).Style = "Normal"

$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = '18.73'
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(27,5).NumberFormat = "@"
$ws.Cells.Item(27,5).Value = '  +11.75%  '
$ws.Cells.Item(27,5).Style = "Normal"

$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value = '1.948'
$ws.Cells.Item(28,4).Style = "Normal"
$ws.Cells.Item(28,5).NumberFormat = "@"
$ws.Cells.Item(28,5).Value = '  +1.92%  '
$ws.Cells.Item(28,5).Style = "Normal"

$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = '1.395'
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Cells.Item(29,5).NumberFormat = "@"
$ws.Cells.Item(29,5).Value = '  +1.44%  '
$ws.Cells.Item(29,5).Style = "Normal"

$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = '4.376'
$ws.Cells.Item(30,4).Style = "Normal"
$ws.Cells.Item(30,5).NumberFormat = "@"
$ws.Cells.Item(30,5).Value = '  +5.33%  '
$ws.Cells.Item(30,5).Style = "Normal"

$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = '0.08975'
$ws.Cells.Item(31,4).Style = "Normal"
$ws.Cells.Item(31,5).NumberFormat = "@"
$ws.Cells.Item(31,5).Value = '  +3.53%  '
$ws.Cells.Item(31,5).Style = "Normal"

$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,4).Value = '4.059'
$ws.Cells.Item(32,4).Style = "Normal"
$ws.Cells.Item(32,5).NumberFormat = "@"
$ws.Cells.Item(32,5).Value = '  +3.11%  '
$ws.Cells.Item(32,5).Style = "Normal"

$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = '0.05303'
$ws.Cells.Item(33,4).Style = "Normal"
$ws.Cells.Item(33,5).NumberFormat = "@"
$ws.Cells.Item(33,5).Value = '  +6.31%  '
$ws.Cells.Item(33,5).Style = "Normal"

$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = '0.7527'
$ws.Cells.Item(34,4).Style = "Normal"
$ws.Cells.Item(34,5).NumberFormat = "@"
$ws.Cells.Item(34,5).Value = '  +6.61%  '
$ws.Cells.Item(34,5).Style = "Normal"

$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value = '1.142'
$ws.Cells.Item(35,4).Style = "Normal"
$ws.Cells.Item(35,5).NumberFormat = "@"
$ws.Cells.Item(35,5).Value = '  +3.33%  '
$ws.Cells.Item(35,5).Style = "Normal"

$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value = '2.692'
$ws.Cells.Item(36,4).Style = "Normal"
$ws.Cells.Item(36,5).NumberFormat = "@"
$ws.Cells.Item(36,5).Value = '  +1.43%  '
$ws.Cells.Item(36,5).Style = "Normal"

$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value = '0.01919'
$ws.Cells.Item(37,4).Style = "Normal"
$ws.Cells.Item(37,5).NumberFormat = "@"
$ws.Cells.Item(37,5).Value = '  +13.52%  '
$ws.Cells.Item(37,5).Style = "Normal"

$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = '2.730'
$ws.Cells.Item(38,4).Style = "Normal"
$ws.Cells.Item(38,5).NumberFormat = "@"
$ws.Cells.Item(38,5).Value = '  +0.65%  '
$ws.Cells.Item(38,5).Style = "Normal"

$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value = '2.183'
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Cells.Item(39,5).NumberFormat = "@"
$ws.Cells.Item(39,5).Value = '  -0.19%  '
$ws.Cells.Item(39,5).Style = "Normal"

$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = '0.9405'
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).NumberFormat = "@"
$ws.Cells.Item(40,5).Value = '  +1.02%  '
$ws.Cells.Item(40,5).Style = "Normal"

$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = '0.4386'
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).NumberFormat = "@"
$ws.Cells.Item(41,5).Value = '  +4.88%  '
$ws.Cells.Item(41,5).Style = "Normal"

$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = '106.09'
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Cells.Item(42,5).NumberFormat = "@"
$ws.Cells.Item(42,5).Value = '  +3.84%  '
$ws.Cells.Item(42,5).Style = "Normal"

$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = '5.880'
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(43,5).NumberFormat = "@"
$ws.Cells.Item(43,5).Value = '  -1.49%  '
$ws.Cells.Item(43,5).Style = "Normal"

$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = '1.001'
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Cells.Item(44,5).NumberFormat = "@"
$ws.Cells.Item(44,5).Value = '  +0.03%  '
$ws.Cells.Item(44,5).Style = "Normal"

$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = '7.747'
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,5).NumberFormat = "@"
$ws.Cells.Item(45,5).Value = '  +4.01%  '
$ws.Cells.Item(45,5).Style = "Normal"

$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = '0.1376'
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Cells.Item(46,5).NumberFormat = "@"
$ws.Cells.Item(46,5).Value = '  +9.49%  '
$ws.Cells.Item(46,5).Style = "Normal"

$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = '0.05846'
$ws.Cells.Item(47,4).Style = "Normal"
$ws.Cells.Item(47,5).NumberFormat = "@"
$ws.Cells.Item(47,5).Value = '  +2.89%  '
$ws.Cells.Item(47,5).Style = "Normal"

$ws.Cells.Item(48,2).NumberFormat = "@"
$ws.Cells.Item(48,2).Value = 'Decentraland'
$ws.Cells.Item(48,2).Style = "Normal"
$ws.Cells.Item(48,3).NumberFormat = "@"
$ws.Cells.Item(48,3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(48,3).Style = "Normal"
$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value = '0.3926'
$ws.Cells.Item(48,4).Style = "Normal"
$ws.Cells.Item(48,5).NumberFormat = "@"
$ws.Cells.Item(48,5).Value = '  +5.94%  '
$ws.Cells.Item(48,5).Style = "Normal"

$ws.Cells.Item(49,2).NumberFormat = "@"
$ws.Cells.Item(49,2).Value = 'EnergySwap'
$ws.Cells.Item(49,2).Style = "Normal"
$ws.Cells.Item(49,3).NumberFormat = "@"
$ws.Cells.Item(49,3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(49,3).Style = "Normal"
$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = '8.604'
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Cells.Item(49,5).NumberFormat = "@"
$ws.Cells.Item(49,5).Value = '  +4.97%  '
$ws.Cells.Item(49,5).Style = "Normal"

$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = '33.58'
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(50,5).NumberFormat = "@"
$ws.Cells.Item(50,5).Value = '  +3.46%  '
$ws.Cells.Item(50,5).Style = "Normal"

$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value = '1.392'
$ws.Cells.Item(51,4).Style = "Normal"
$ws.Cells.Item(51,5).NumberFormat = "@"
$ws.Cells.Item(51,5).Value = '  +4.02%  '
$ws.Cells.Item(51,5).Style = "Normal"
